$d = $word.ActiveDocument
$t = $d.Tables(1)

# New lattice-multiplication problems, in row-major order (5 rows x 3 cols).
# Each entry: line1 "A x B", line2 "  d1    d2" (digits of B), line3 "  ----",
# line4 "d1|    |", line5 "d2|    |" (digits of A).
$cells = @(
  @("61 x 32", "  3    2", "  ----", "6|    |", "1|    |"),
  @("73 x 45", "  4    5", "  ----", "7|    |", "3|    |"),
  @("81 x 27", "  2    7", "  ----", "8|    |", "1|    |"),
  @("79 x 42", "  4    2", "  ----", "7|    |", "9|    |"),
  @("64 x 91", "  9    1", "  ----", "6|    |", "4|    |"),
  @("94 x 87", "  8    7", "  ----", "9|    |", "4|    |"),
  @("28 x 26", "  2    6", "  ----", "2|    |", "8|    |"),
  @("95 x 76", "  7    6", "  ----", "9|    |", "5|    |"),
  @("97 x 85", "  8    5", "  ----", "9|    |", "7|    |"),
  @("69 x 18", "  1    8", "  ----", "6|    |", "9|    |"),
  @("13 x 96", "  9    6", "  ----", "1|    |", "3|    |"),
  @("28 x 36", "  3    6", "  ----", "2|    |", "8|    |"),
  @("55 x 97", "  9    7", "  ----", "5|    |", "5|    |"),
  @("80 x 73", "  7    3", "  ----", "8|    |", "0|    |"),
  @("19 x 43", "  4    3", "  ----", "1|    |", "9|    |")
)

$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$idx = 0
for ($row = 1; $row -le 5; $row++) {
  for ($col = 1; $col -le 3; $col++) {
    $lines = $cells[$idx]
    $l1 = $lines[0]
    $l2 = $lines[1]
    $l3 = $lines[2]
    $l4 = $lines[3]
    $l5 = $lines[4]

    # Lines 2 and 3 start with literal spaces, so they need xml:space="preserve"
    # to round-trip exactly like the source document.
    $frag = '<w:p ' + $wns + '><w:r><w:rPr><w:sz w:val="32"/></w:rPr>' + `
      '<w:t>' + $l1 + '</w:t><w:br/>' + `
      '<w:t xml:space="preserve">' + $l2 + '</w:t><w:br/>' + `
      '<w:t xml:space="preserve">' + $l3 + '</w:t><w:br/>' + `
      '<w:t>' + $l4 + '</w:t><w:br/>' + `
      '<w:t>' + $l5 + '</w:t>' + `
      '</w:r></w:p>'

    $cell = $t.Cell($row, $col)
    $cellRange = $cell.Range
    $full = $cellRange.Text
    # Exclude the trailing paragraph mark + cell-end mark (2 chars) so the
    # inserted paragraph replaces the cell's content in place.
    $target = $d.Range($cellRange.Start, $cellRange.Start + $full.Length - 2)
    $target.InsertXML($frag)

    $idx++
  }
}
